# edit.ps1 -- PowerShell-style PowerPoint COM-interop script
#
# Reproduces the two changes described by the commit diff:
#
#   1. The table on slide 5 (the "B1- TYPES OF FINANCIAL DOCUMENTS" table)
#      switches from the deck's custom table style
#      {ED1342F1-C580-4AF8-B829-0FFC22B27A27} to the built-in table style
#      {6B5BE778-6056-4731-9F5D-80733107F8D3}.
#
#   2. The presentation's design colors change from the custom "Red Violet"
#      palette to the standard Office palette (dk1/lt1/dk2/lt2/accent1-6/
#      hlink/folHlink), i.e. the Colors swatch on the Design tab is switched
#      from the deck's bespoke "Integral" colors to plain "Office" colors.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 5
# ---------------------------------------------------------------------

$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6B5BE778-6056-4731-9F5D-80733107F8D3}")
    }
}

# ---------------------------------------------------------------------
# 2) Theme colors: Red Violet -> Office
# ---------------------------------------------------------------------

function ConvertTo-RgbValue {
    param([string]$HexColor)
    $r = [Convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colors, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-RgbValue $officeColors[$i - 1]
}
